$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns E and F (header labels "soft"/"rigid" and all their values)
# for every used row (1 through 21) on the sheet.
$lastRow = 21

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)  # column E
    $fCell = $ws.Cells.Item($r, 6)  # column F

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    $eCell.Value = $fVal
    $fCell.Value = $eVal
}
